$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows before the current row 4 (old rows 4-7 shift down to 6-9)
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# --- New row 4 (brand-new weekly record) ---
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value = 44435
$ws.Range("D4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100108
$ws.Range("H4").Value = "Tropicales y subtropicales"
$ws.Range("I4").Value = 100108001
$ws.Range("J4").Value = "Guayaba"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 130
$ws.Range("N4").Value = 1300
$ws.Range("O4").Value = 1300
$ws.Range("P4").Value = 1300
$ws.Range("Q4").Value = "$/kilo"
$ws.Range("R4").Value = "Región de Arica y Parinacota"
$ws.Range("S4").Value = 1300
$ws.Range("T4").Value = 1

# --- New row 5 (brand-new weekly record) ---
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 44431
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100108
$ws.Range("H5").Value = "Tropicales y subtropicales"
$ws.Range("I5").Value = 100108001
$ws.Range("J5").Value = "Guayaba"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 1300
$ws.Range("O5").Value = 1300
$ws.Range("P5").Value = 1300
$ws.Range("Q5").Value = "$/kilo"
$ws.Range("R5").Value = "Región de Arica y Parinacota"
$ws.Range("S5").Value = 1300
$ws.Range("T5").Value = 1

# --- Append new row 10 at the end (brand-new weekly record) ---
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = 44432
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100108
$ws.Range("H10").Value = "Tropicales y subtropicales"
$ws.Range("I10").Value = 100108001
$ws.Range("J10").Value = "Guayaba"
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 1300
$ws.Range("O10").Value = 1300
$ws.Range("P10").Value = 1300
$ws.Range("Q10").Value = "$/kilo"
$ws.Range("R10").Value = "Región de Arica y Parinacota"
$ws.Range("S10").Value = 1300
$ws.Range("T10").Value = 1
